# Generate Report for Handback
# - Update the Status text for the 43a8a520 row (row 3) in both the
#   zh-cn and de-de sheets from "Ready for handoff" to
#   "Handback transform failed".
# - Fill in the "Error Detail" cell (column P, row 3) on both the
#   zh-cn and de-de sheets with a message describing the handback/
#   handoff file-name mismatch.
# - Widen column P on both sheets to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcnError = "Handback file name: 5gqifsla.35o is different with handoff file name: 43a8a520-5615-49ca-941b-aa362aa96267.ee4bcaf9eba7f635377f76beb64e4ecc46a7324f.zh-cn."
$dedeError = "Handback file name: 5gqifsla.35o is different with handoff file name: 43a8a520-5615-49ca-941b-aa362aa96267.ee4bcaf9eba7f635377f76beb64e4ecc46a7324f.de-de."

$zhcn.Range("P3").Value = $zhcnError
$dede.Range("P3").Value = $dedeError

# Excel stores column widths in a character-width unit that is rounded to
# whole pixels internally; asking for 39.17 is what round-trips to an
# on-disk <col> width of exactly 40 (same quirk already visible on column A).
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
